# Auto-generated Excel COM-interop script
# Applies scheduled-runner value updates across ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 17: One for the Road
$ws.Range("H17").Value = 1624961.1
$ws.Range("J17").Value = 1670093.4
$ws.Range("L17").Value = 5010280.199999999
$ws.Range("N17").Value = -5010616.199999999

# Row 28: The Writing Is Not on the Wall
$ws.Range("H28").Value = 794.4286
$ws.Range("J28").Value = 1099.5
$ws.Range("L28").Value = 1099.5
$ws.Range("N28").Value = -2069.5

# Row 88: The Grave of Hemlock Groves
$ws.Range("H88").Value = 1458
$ws.Range("I88").Value = 1460.3334
$ws.Range("J88").Value = 1455.2
$ws.Range("K88").Value = 1460.3334
$ws.Range("L88").Value = 1455.2
$ws.Range("M88").Value = -1054.3334
$ws.Range("N88").Value = -2267.2

# Row 91: Dappling the Highlands (L)
$ws.Range("H91").Value = 1458
$ws.Range("I91").Value = 1460.3334
$ws.Range("J91").Value = 1455.2
$ws.Range("K91").Value = 1460.3334
$ws.Range("L91").Value = 1455.2
$ws.Range("M91").Value = -56.33339999999998
$ws.Range("N91").Value = -4263.2

# Row 137: Cutting Edge of Culinary Quality
$ws.Range("H137").Value = 257776.22
$ws.Range("I137").Value = 19990.5
$ws.Range("J137").Value = 325715
$ws.Range("K137").Value = 59971.5
$ws.Range("L137").Value = 977145
$ws.Range("M137").Value = -57421.5
$ws.Range("N137").Value = -982245

# Row 141: Remedy for Reason
$ws.Range("H141").Value = 2941.1667
$ws.Range("I141").Value = 2723.0667
$ws.Range("K141").Value = 8169.2001
$ws.Range("M141").Value = -2989.2001

$ws = $wb.Worksheets.Item("ARM")
# Row 32: Ingot We Trust
$ws.Range("H32").Value = 2861500
$ws.Range("I32").Value = 3335716.8
$ws.Range("K32").Value = 3335716.8
$ws.Range("M32").Value = -3335429.8

# Row 61: Dealing with the Tough Stuff
$ws.Range("H61").Value = 926633.0600000001
$ws.Range("I61").Value = 4302.893
$ws.Range("K61").Value = 4302.893
$ws.Range("M61").Value = -4090.893

# Row 132: Don't Bore Me, Ore Me
$ws.Range("H132").Value = 4141141.5
$ws.Range("I132").Value = 1709.1923
$ws.Range("J132").Value = 22078682
$ws.Range("K132").Value = 5127.5769
$ws.Range("L132").Value = 66236046
$ws.Range("M132").Value = -2597.5769
$ws.Range("N132").Value = -66241106

# Row 136: Metal with Mettle
$ws.Range("H136").Value = 926633.0600000001
$ws.Range("I136").Value = 4302.893
$ws.Range("K136").Value = 12908.679
$ws.Range("M136").Value = -10358.679

$ws = $wb.Worksheets.Item("BSM")
# Row 134: Ruthenium Supremium
$ws.Range("H134").Value = 17815.777
$ws.Range("I134").Value = 6973.5386
$ws.Range("K134").Value = 20920.6158
$ws.Range("M134").Value = -18385.6158

# Row 139: Maul Me
$ws.Range("H139").Value = 267572
$ws.Range("J139").Value = 267572
$ws.Range("L139").Value = 267572
$ws.Range("N139").Value = -277852

$ws = $wb.Worksheets.Item("CRP")
# Row 22: Driving Up the Wall
$ws.Range("H22").Value = 1481.5883
$ws.Range("I22").Value = 1181.4286
$ws.Range("K22").Value = 1181.4286
$ws.Range("M22").Value = -831.4286

# Row 31: Wall Not Found
$ws.Range("H31").Value = 287503
$ws.Range("I31").Value = 514999
$ws.Range("K31").Value = 514999
$ws.Range("M31").Value = -514704

# Row 34: Armoires of the Rich and Famous
$ws.Range("H34").Value = 287503
$ws.Range("I34").Value = 514999
$ws.Range("K34").Value = 514999
$ws.Range("M34").Value = -514797

# Row 105: Zelkova, My Love
$ws.Range("H105").Value = 12412.223
$ws.Range("I105").Value = 15387.143
$ws.Range("K105").Value = 15387.143
$ws.Range("M105").Value = -13640.143

# Row 132: Hull Lotta Damage
$ws.Range("H132").Value = 41178976
$ws.Range("I132").Value = 2349.158
$ws.Range("K132").Value = 7047.474
$ws.Range("M132").Value = -4517.474

# Row 133: Yimepi's Country Charms
$ws.Range("H133").Value = 77116.94
$ws.Range("J133").Value = 77116.94
$ws.Range("L133").Value = 77116.94
$ws.Range("N133").Value = -82176.94

$ws = $wb.Worksheets.Item("CUL")
# Row 14: Keep Your Powder Dry
$ws.Range("H14").Value = 0
$ws.Range("I14").Value = 0
$ws.Range("K14").Value = 0
$ws.Range("M14").ClearContents()

# Row 41: Gegeruju Gets Down
$ws.Range("H41").Value = 254070.83
$ws.Range("I41").Value = 775
$ws.Range("J41").Value = 304730
$ws.Range("K41").Value = 2325
$ws.Range("L41").Value = 914190
$ws.Range("M41").Value = -1987
$ws.Range("N41").Value = -914866

# Row 56: Culture Club
$ws.Range("H56").Value = 6898.3335
$ws.Range("I56").Value = 6898.3335
$ws.Range("K56").Value = 6898.3335
$ws.Range("M56").Value = -6368.3335

# Row 107: Slippery Service
$ws.Range("H107").Value = 1183.1666
$ws.Range("J107").Value = 1500
$ws.Range("L107").Value = 4500
$ws.Range("N107").Value = -8340

$ws = $wb.Worksheets.Item("GSM")
# Row 122: Awarding Academic Excellence
$ws.Range("H122").Value = 3657
$ws.Range("I122").Value = 2516.5
$ws.Range("K122").Value = 7549.5
$ws.Range("M122").Value = -5099.5

# Row 132: On Board for Lar
$ws.Range("H132").Value = 1552516.9
$ws.Range("I132").Value = 3168.889
$ws.Range("K132").Value = 9506.667000000001
$ws.Range("M132").Value = -6976.667000000001

$ws = $wb.Worksheets.Item("LTW")
# Row 46: Supply Side Logic
$ws.Range("H46").Value = 1612.25
$ws.Range("I46").Value = 1300
$ws.Range("J46").Value = 1799.6
$ws.Range("K46").Value = 1300
$ws.Range("L46").Value = 1799.6
$ws.Range("M46").Value = -1112
$ws.Range("N46").Value = -2175.6

# Row 53: Foot Blues
$ws.Range("H53").Value = 7025
$ws.Range("I53").Value = 7000
$ws.Range("J53").Value = 7050
$ws.Range("K53").Value = 7000
$ws.Range("L53").Value = 7050
$ws.Range("M53").Value = -6482
$ws.Range("N53").Value = -8086

# Row 55: It's Not a Job, It's a Calling
$ws.Range("H55").Value = 1779.2609
$ws.Range("I55").Value = 2041.25
$ws.Range("J55").Value = 1639.5333
$ws.Range("K55").Value = 2041.25
$ws.Range("L55").Value = 1639.5333
$ws.Range("M55").Value = -1868.25
$ws.Range("N55").Value = -1985.5333

# Row 68: You Could Say It's a Moving Target
$ws.Range("H68").Value = 7688.0557
$ws.Range("I68").Value = 9752.083000000001
$ws.Range("J68").Value = 3560
$ws.Range("K68").Value = 9752.083000000001
$ws.Range("L68").Value = 3560
$ws.Range("M68").Value = -9003.083000000001
$ws.Range("N68").Value = -5058

# Row 71: They Call It Bloody Mary (L)
$ws.Range("H71").Value = 7688.0557
$ws.Range("I71").Value = 9752.083000000001
$ws.Range("J71").Value = 3560
$ws.Range("K71").Value = 48760.415
$ws.Range("L71").Value = 17800
$ws.Range("M71").Value = -45016.415
$ws.Range("N71").Value = -25288

# Row 82: Trainin' the Neck
$ws.Range("H82").Value = 2475.5
$ws.Range("J82").Value = 2333.3333
$ws.Range("L82").Value = 2333.3333
$ws.Range("N82").Value = -3055.3333

# Row 85: Training Is Only Skintight (L)
$ws.Range("H85").Value = 2475.5
$ws.Range("J85").Value = 2333.3333
$ws.Range("L85").Value = 2333.3333
$ws.Range("N85").Value = -4829.3333

# Row 139: Giving Gatherers Their Gear
$ws.Range("H139").Value = 94973.75
$ws.Range("J139").Value = 94973.75
$ws.Range("L139").Value = 94973.75
$ws.Range("N139").Value = -105253.75

$ws = $wb.Worksheets.Item("WVR")
# Row 21: Don't Trew So Hard
$ws.Range("H21").Value = 53750
$ws.Range("I21").Value = 30000
$ws.Range("K21").Value = 30000
$ws.Range("M21").Value = -29765

# Row 35: Pantser Corps
$ws.Range("H35").Value = 53750
$ws.Range("I35").Value = 30000
$ws.Range("K35").Value = 30000
$ws.Range("M35").Value = -29710

# Row 62: Pride Up in Smoke
$ws.Range("H62").Value = 23728.6
$ws.Range("I62").Value = 14124.143
$ws.Range("K62").Value = 14124.143
$ws.Range("M62").Value = -13500.143

# Row 65: Desperate for Diversionaries (L)
$ws.Range("H65").Value = 23728.6
$ws.Range("I65").Value = 14124.143
$ws.Range("K65").Value = 70620.715
$ws.Range("M65").Value = -67500.715

# Row 107: Flax Wax
$ws.Range("H107").Value = 856.5
$ws.Range("I107").Value = 939.4737
$ws.Range("J107").Value = 541.2
$ws.Range("K107").Value = 2818.4211
$ws.Range("L107").Value = 1623.6
$ws.Range("M107").Value = -898.4211
$ws.Range("N107").Value = -5463.6

# Row 132: Comfy Cabins
$ws.Range("H132").Value = 839351.5600000001
$ws.Range("I132").Value = 2201
$ws.Range("K132").Value = 6603
$ws.Range("M132").Value = -4073

# Row 136: Weaving the Envelope
$ws.Range("H136").Value = 304925.88
$ws.Range("I136").Value = 2357.0476
$ws.Range("K136").Value = 7071.1428
$ws.Range("M136").Value = -4521.1428
